$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 86 (shifts old rows 86-92 down to 87-93)
$ws.Rows("86:86").Insert()

# Fill newly inserted row 86 with new case data
$ws.Range("A86").Value = "'7237"
$ws.Range("B86").Value = "'9/3/2025"
$ws.Range("C86").Value = "NAZCA 3985"
$ws.Range("D86").Value = "'15"
$ws.Range("E86").Value = "'809409491"
$ws.Range("F86").Value = "AYKO"
$ws.Range("G86").Value = "Pendiente"
$ws.Range("H86").Value = "Picada"
$ws.Range("I86").Value = 1
$ws.Range("J86").Value = "Cambio"
$ws.Range("K86").Value = "Sin equipos"
$ws.Range("L86").Value = "Pasante"
$ws.Range("M86").Value = -58.497366
$ws.Range("N86").Value = -34.591544
$ws.Range("O86").Value = "Paternal"
$ws.Range("P86").Value = "Capital Norte"

# Update fields on row 90 (was old row 89, shifted by insert): OT -> Pendiente ADM, Observaciones changed
$ws.Range("E90").Value = "Pendiente ADM"
$ws.Range("H90").Value = "Desmonte de columna"

# Update fields on row 91 (was old row 90, shifted by insert): Caso, Direccion and coordinates changed
$ws.Range("A91").Value = "'7247"
$ws.Range("C91").Value = "ARAUJO 3430"
$ws.Range("M91").Value = -58.46866
$ws.Range("N91").Value = -34.669465

# Append new row 94 (Caso 7240)
$ws.Range("A94").Value = "'7240"
$ws.Range("B94").Value = "'9/16/2025"
$ws.Range("C94").Value = "LARRALDE, CRISOLOGO AV. 3875"
$ws.Range("D94").Value = "'12"
$ws.Range("E94").Value = "'809784524"
$ws.Range("F94").Value = "AYKO"
$ws.Range("G94").Value = "Pendiente"
$ws.Range("H94").Value = "Picada"
$ws.Range("I94").Value = 1
$ws.Range("J94").Value = "Cambio"
$ws.Range("K94").Value = "Sin equipos"
$ws.Range("L94").Value = "Terminal"
$ws.Range("M94").Value = -58.481316
$ws.Range("N94").Value = -34.556157
$ws.Range("O94").Value = "Saavedra"
$ws.Range("P94").Value = "Capital Norte"

# Append new row 95 (Caso 7248)
$ws.Range("A95").Value = "'7248"
$ws.Range("B95").Value = "'9/16/2025"
$ws.Range("C95").Value = "FERNANDEZ DE LA CRUZ, F., GRAL. AV. 4065"
$ws.Range("D95").Value = "'8"
$ws.Range("E95").Value = "'809784526"
$ws.Range("F95").Value = "AYKO"
$ws.Range("G95").Value = "Pendiente"
$ws.Range("H95").Value = "Columna chocada "
$ws.Range("I95").Value = 1
$ws.Range("J95").Value = "Cambio"
$ws.Range("K95").Value = "Sin equipos"
$ws.Range("L95").Value = "Pasante"
$ws.Range("M95").Value = -58.455155
$ws.Range("N95").Value = -34.669378
$ws.Range("O95").Value = "Boedo"
$ws.Range("P95").Value = "Capital Sur"
